$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.901.62'
$ws.Range('E2').Value = '  +0.28%  '

$ws.Range('D3').Value = '1.884.94'
$ws.Range('E3').Value = '  -0.11%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '336.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.58%  '

$ws.Range('E6').Value = '  +0.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4703'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.50%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3969'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.85%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08057'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.33%  '

$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.63'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.36%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.018'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.84%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.02%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.879.44'
$ws.Range('E13').Value = '  -0.34%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.031'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.59%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.312'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.34%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.011'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.09%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.32'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.29%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06735'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.21%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001050'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.12%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.20%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.009'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.21%  '

$ws.Range('D22').Value = '27.877.50'
$ws.Range('E22').Value = '  +0.11%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.516'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.58%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.316'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.34%  '

$ws.Range('D26').Value = '2.090.97'
$ws.Range('E26').Value = '  -0.90%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.07%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.90'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.43%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.167'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.74%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.522'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.28%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '122.21'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.11%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9896'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.77%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09528'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.43%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.638'
$ws.Range('D34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.361'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.361'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.32%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06102'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.98%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02250'
$ws.Range('D38').Style = 'Normal'

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.342'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.85%  '

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.201'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.55%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.008'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6014'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.01%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1903'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.21%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.05%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5701'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.32%  '

$ws.Range('E46').Value = '  -0.87%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.39%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.951'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.14%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06796'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.72%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '112.86'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.99%  '

$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.070'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.22%  '
